$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E9").Value = "Não iniciado"
